# Outstandings.xlsx - "Add files via upload"
#
# Adds a new purchase-outstanding line item (invoice "08/23-24", dated
# 08/23/2023, for Namrata Rubber Product Pvt Ltd, amount 29140) to the
# first group of the "Purchase 22-23" sheet, updates the group's running
# total formula, and renumbers/shifts the rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8; everything at/after row 8 (the second
# and third groups) shifts down by one row, and Excel auto-adjusts the
# existing formulas (F12/F14) that reference those shifted rows.
$ws.Rows("8:8").Insert()

# Give the new row 8 the same visual formatting (fonts/borders/number
# formats/alignment + row height) as the other detail rows in this
# group (row 2 is an identical-looking "group" row: Sr.No/Date/Invoice/
# Vendor/Outstanding columns with the thin-border boxed style).
$ws.Range("A2:F2").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122) | Out-Null
$ws.Rows("8:8").RowHeight = 14.4

# Fill in the new line item's data.
$ws.Range("B8").Value = 45048
$ws.Range("C8").Value = "08/23-24"
$ws.Range("D8").Value = "Namrata Rubber Product Pvt Ltd"
$ws.Range("E8").Value = 29140

# The group subtotal formula moves from F7 down to F8, now summing the
# five rows of the group (E4:E8).
$ws.Range("F7").ClearContents()
$ws.Range("F8").Formula = "=E4+E5+E6+E7+E8"

# Restore the sheet's saved selection/cursor position.
$ws.Range("A16").Select() | Out-Null
